$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "578.99" or
# "65.005.88"); force text entry via the "@" number format so Excel
# keeps the exact display string, then restore the default style so
# the cell is not left with a lingering custom format.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '65.005.88'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +1.88%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.175.82'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +4.03%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '578.99'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +3.46%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '152.09'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +6.68%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.174.97'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +4.06%  '
$ws.Range('E9').Value = '  +3.68%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.162'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +5.93%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.24'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.51%  '
$ws.Range('E12').Value = '  +2.83%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000271'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +18.00%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '37.88'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +6.06%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.696.98'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +4.02%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.104.76'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.89%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.182.18'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +4.20%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '7.18'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +5.55%  '
$ws.Range('E19').Value = '  +1.32%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '514.91'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +8.21%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.88'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +5.83%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.731'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +6.91%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '15.30'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +3.83%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '7.81'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +3.62%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '85.46'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +3.26%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.93'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +4.81%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.98'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +10.21%  '
$ws.Range('E29').Value = '  +7.35%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '27.94'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +6.33%  '
$ws.Range('E31').Value = '  +13.55%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.20'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +5.14%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.34'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +9.79%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.58'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +5.71%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '55.81'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.25%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0899'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +10.42%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.16'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +12.64%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '475.41'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +6.51%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0422'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +3.13%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '8.67'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +4.76%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.071.11'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.53%  '
$ws.Range('E43').Value = '  +1.66%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.286'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +5.98%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.40'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +6.23%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '29.12'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +2.90%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0₃0612'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +19.10%  '
$ws.Range('E49').Value = '  +2.35%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '120.51'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.14%  '
